# "updated fan link in bom"
#
# Component Part List:
#  - B11: "GT2 Idler Pulley 20NT"   -> "GT2 Idler Pulley Smooth"
#  - B19: "Nema 17 Stepper Motor"   -> "1.7A Nema 17 Stepper Motor"   (D19 qty 5 -> 4)
#  - D25: qty 5 -> 4 (DRV 8825 / ZE Stepper Drivers)
#  - B28: "12V Layer Fan"           -> "12/24V Layer Fan"             (E28 price 10 -> 5)
#  - Row 30 ("24V to 12V DC Converter" / "1.5A 18W DC DC Converter") removed entirely,
#    along with its hyperlink. Totals (E30, formerly E31) recalc automatically.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Component Part List")

# Rename the GT2 idler pulley variant to the "Smooth" part.
$ws.Range("B11").Value = "GT2 Idler Pulley Smooth"

# Nema 17 stepper motor -> higher current variant, quantity drops from 5 to 4.
$ws.Range("B19").Value = "1.7A Nema 17 Stepper Motor"
$ws.Range("D19").Value = 4

# ZE stepper driver (DRV 8825) order quantity drops from 5 to 4.
$ws.Range("D25").Value = 4

# Layer fan is now rated for both 12V and 24V, and its bulk price dropped to 5.
$ws.Range("B28").Value = "12/24V Layer Fan"
$ws.Range("E28").Value = 5

# Remove the obsolete "24V to 12V DC Converter" row (row 30) and its hyperlink.
foreach ($hl in $ws.Hyperlinks) {
    if ($hl.Range.Address() -eq "`$B`$30") {
        $hl.Delete()
    }
}
$ws.Rows("30:30").Delete()

"done"
